$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '34.246.37'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '1.829.10'
$ws.Range('E3').Value = '  +2.73%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '225.28'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').Value = '0.558'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '32.11'
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('D9').Value = '0.291'
$ws.Range('E9').Value = '  +3.88%  '
$ws.Range('D10').Value = '0.0720'
$ws.Range('E10').Value = '  +9.54%  '
$ws.Range('D11').Value = '0.0930'
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').Value = '2.091.61'
$ws.Range('E12').Value = '  +2.76%  '
$ws.Range('D13').Value = '1.829.59'
$ws.Range('E13').Value = '  +2.77%  '
$ws.Range('D14').Value = '10.86'
$ws.Range('E14').Value = '  -2.73%  '
$ws.Range('D15').Value = '0.644'
$ws.Range('E15').Value = '  +2.88%  '
$ws.Range('D16').Value = '34.238.77'
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').Value = '4.34'
$ws.Range('E17').Value = '  +3.02%  '
$ws.Range('D18').Value = '69.75'
$ws.Range('E18').Value = '  +1.45%  '
$ws.Range('D19').Value = '252.24'
$ws.Range('E19').Value = '  -1.17%  '
$ws.Range('E20').Value = '  +6.82%  '
$ws.Range('D21').Value = '11.22'
$ws.Range('E21').Value = '  +8.25%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = '4.28'
$ws.Range('E23').Value = '  +1.99%  '
$ws.Range('D24').Value = '2.17'
$ws.Range('E24').Value = '  +1.40%  '
$ws.Range('D25').Value = '160.55'
$ws.Range('E25').Value = '  +2.14%  '
$ws.Range('D26').Value = '16.69'
$ws.Range('E26').Value = '  +1.68%  '
$ws.Range('D27').Value = '7.27'
$ws.Range('E27').Value = '  +3.73%  '
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('B30').Value = 'WavesCommunityToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/U31RirsudBxis+wavescommunitytoken-wct'
$ws.Range('D30').Value = '197.54'
$ws.Range('E30').Value = '  +6.79 million%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.0537'
$ws.Range('E31').Value = '  +4.72%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '3.79'
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.21'
$ws.Range('E33').Value = '  +2.10%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '3.59'
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '1.90'
$ws.Range('E35').Value = '  +1.97%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '1.443.90'
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '0.646'
$ws.Range('E37').Value = '  +3.54%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '1.07'
$ws.Range('E38').Value = '  +1.54%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.0190'
$ws.Range('E39').Value = '  +1.68%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '0.968'
$ws.Range('E40').Value = '  +8.66%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '81.90'
$ws.Range('E41').Value = '  -1.24%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = '2.76'
$ws.Range('E42').Value = '  -3.18%  '
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').Value = '2.35'
$ws.Range('E43').Value = '  +0.45%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '2.16'
$ws.Range('E44').Value = '  +5.18%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '6.09'
$ws.Range('E45').Value = '  +4.48%  '
$ws.Range('D46').Value = '1.987.65'
$ws.Range('E46').Value = '  +2.53%  '
$ws.Range('B47').Value = 'Kaspa'
$ws.Range('C47').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D47').Value = '0.0497'
$ws.Range('E47').Value = '  -2.90%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '1.05'
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '106.57'
$ws.Range('E49').Value = '  +8.25%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '11.89'
$ws.Range('E51').Value = '  -2.92%  '
